$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewOR")

# Update Status for newInsOR_001 from 'E' to 'P'
$ws.Range("B2").Value = "P"

# Update DepVal for newInsOR_001 to the new pending invoice id
$ws.Range("D2").Value = "EAOR21AP-0316"
